# Update the WRESBAL FRED data workbook:
#  - Append a new weekly observation row to the "Data" sheet.
#  - Refresh the FRED series metadata on the "SeriesInfo" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Data" sheet: append row 112 (new weekly WRESBAL observation)
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Copy the formatting (date number format, borders, alignment) of the last
# existing data row down onto the new row before filling in values, so the
# new date cell matches the style used by every other date cell in column A.
$wsData.Range("A111").Copy()
$wsData.Range("A112").PasteSpecial(-4122)  # xlPasteFormats

$wsData.Range("A112").Value = 45245
$wsData.Range("B112").Value = 3391.713

# ---------------------------------------------------------------------------
# "SeriesInfo" sheet: refresh metadata pulled from the FRED API
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# realtime_start / realtime_end / observation_end hold plain textual dates
# (not real Excel dates) in the source data. A leading apostrophe keeps
# them stored as text instead of letting them be auto-recognized as dates.
$wsInfo.Range("B3").Value = "'2023-11-21"
$wsInfo.Range("B4").Value = "'2023-11-21"
$wsInfo.Range("B7").Value = "'2023-11-15"

$wsInfo.Range("B14").Value = "2023-11-16 15:36:02-06"
$wsInfo.Range("B15").Value = 74
